# Insert a new data row at row 456 (pushing all existing rows 456..533 down
# to 457..534), then populate the newly inserted row with the new record's
# values. This mirrors the source diff, where dimension grows from
# A1:R533 to A1:R534 and every row from 457..534 carries what used to be
# the row above it, with row 456 holding brand-new data
# (Paine / 1a (guarda) / 2023-11-28).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("456:456").Insert()

$ws.Cells.Item(456, 1).Value = 5
$ws.Cells.Item(456, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(456, 3).Value = "Maule"
$ws.Cells.Item(456, 4).Value = 45258
$ws.Cells.Item(456, 5).Value = 7
$ws.Cells.Item(456, 6).Value = 100112045
$ws.Cells.Item(456, 7).Value = "Zapallo"
$ws.Cells.Item(456, 8).Value = "Paine"
$ws.Cells.Item(456, 9).Value = "1a (guarda)"
$ws.Cells.Item(456, 10).Value = 1000
$ws.Cells.Item(456, 11).Value = 900
$ws.Cells.Item(456, 12).Value = 900
$ws.Cells.Item(456, 13).Value = 900
$ws.Cells.Item(456, 14).Value = "`$/kilo (volumen en unidades)"
$ws.Cells.Item(456, 15).Value = "Región del Maule"
$ws.Cells.Item(456, 16).Value = 900
$ws.Cells.Item(456, 17).Value = 1
$ws.Cells.Item(456, 18).Value = "Hortaliza"
